# Update "想去人数" (F column) counts on both the "展览" and "全部类型" sheets.
$wb = $excel.ActiveWorkbook

$updates = @{
    2  = 1137
    3  = 847
    4  = 278
    8  = 2373
    9  = 7743
    10 = 919
    12 = 382
    14 = 427
    17 = 7947
    19 = 1381
    20 = 157
    24 = 324
    25 = 164
    29 = 26
    30 = 423
    37 = 80
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
